# Generate Report for Handback
# Updates the localization-status workbook to reflect that the handback
# has completed and is in sync with en-US: status text, handback
# timestamps, clearing of the stale "handback not latest" error details,
# and widening of the Status / Error Detail columns to fit the new text.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

# Column widths (as exposed through the ColumnWidth COM property) that
# reproduce the wider "Status"/"Error Detail" columns used to show the
# longer handback message.
$wideStatusWidth = 29.166666666666668
$wideErrorWidth  = 12.833333333333332

# --- Overview sheet ---------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("E2").Value = $newStatus
$ov.Range("F2").Value = $newStatus
$ov.Range("E3").Value = $newStatus
$ov.Range("F3").Value = $newStatus

$ov.Columns.Item(5).ColumnWidth = $wideStatusWidth
$ov.Columns.Item(6).ColumnWidth = $wideStatusWidth

# --- zh-cn sheet --------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("C2").Value = $newStatus
$zh.Range("C3").Value = $newStatus

$zh.Range("J2").Value = "2016-08-02 09:42:21"
$zh.Range("J3").Value = "2016-08-02 09:42:21"

$zh.Range("O2").Value = ""
$zh.Range("O3").Value = ""

$zh.Columns.Item(3).ColumnWidth = $wideStatusWidth
$zh.Columns.Item(15).ColumnWidth = $wideErrorWidth

# --- de-de sheet --------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("C2").Value = $newStatus
$de.Range("C3").Value = $newStatus

$de.Range("J2").Value = "2016-08-02 09:42:36"
$de.Range("J3").Value = "2016-08-02 09:42:36"

$de.Range("O2").Value = ""
$de.Range("O3").Value = ""

$de.Columns.Item(3).ColumnWidth = $wideStatusWidth
$de.Columns.Item(15).ColumnWidth = $wideErrorWidth
